$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for rows 2-14 from serial 45180 to 45181
for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 3).Value = 45181
}
